$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.214.57"
$ws.Range("D3").Value = "2.275.07"
$ws.Range("E3").Value = "  -2.77%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'298.31"
$ws.Range("E5").Value = "  -2.74%  "
$ws.Range("D6").Value = "'95.27"
$ws.Range("E6").Value = "  -5.85%  "
$ws.Range("E8").Value = "  -3.91%  "
$ws.Range("E9").Value = "  -3.83%  "
$ws.Range("D10").Value = "'33.32"
$ws.Range("E10").Value = "  -4.56%  "
$ws.Range("E11").Value = "  -1.18%  "
$ws.Range("D12").Value = "'48.14"
$ws.Range("E12").Value = "  -8.30%  "
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").Value = "'6.64"
$ws.Range("E14").Value = "  -3.57%  "
$ws.Range("D15").Value = "'15.75"
$ws.Range("E15").Value = "  -0.72%  "
$ws.Range("D16").Value = "2.628.45"
$ws.Range("E16").Value = "  -2.80%  "
$ws.Range("D17").Value = "2.267.26"
$ws.Range("E17").Value = "  -4.96%  "
$ws.Range("D18").Value = "'0.780"
$ws.Range("E18").Value = "  -5.86%  "
$ws.Range("D19").Value = "42.196.06"
$ws.Range("E20").Value = "  -2.39%  "
$ws.Range("D21").Value = "'11.46"
$ws.Range("E21").Value = "  -2.60%  "
$ws.Range("E22").Value = "  -3.94%  "
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("D24").Value = "'233.08"
$ws.Range("E24").Value = "  -1.57%  "
$ws.Range("D25").Value = "'1.97"
$ws.Range("E25").Value = "  -3.12%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("E27").Value = "  -4.57%  "
$ws.Range("D28").Value = "'23.87"
$ws.Range("E28").Value = "  -6.20%  "
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("D30").Value = "'166.86"
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("D31").Value = "'33.92"
$ws.Range("E31").Value = "  -4.16%  "
$ws.Range("D32").Value = "'9.05"
$ws.Range("E32").Value = "  -3.17%  "
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D34").Value = "'4.92"
$ws.Range("E34").Value = "  -4.10%  "
$ws.Range("D35").Value = "'4.50"
$ws.Range("E35").Value = "  -3.12%  "
$ws.Range("D36").Value = "'0.0691"
$ws.Range("E36").Value = "  -5.23%  "
$ws.Range("E37").Value = "  -5.01%  "
$ws.Range("D38").Value = "'16.17"
$ws.Range("E38").Value = "  -8.20%  "
$ws.Range("D39").Value = "'2.79"
$ws.Range("E39").Value = "  -4.47%  "
$ws.Range("D40").Value = "'0.0989"
$ws.Range("E40").Value = "  -3.08%  "
$ws.Range("E41").Value = "  -3.51%  "
$ws.Range("D42").Value = "'1.73"
$ws.Range("E42").Value = "  -7.27%  "
$ws.Range("E43").Value = "  -4.47%  "
$ws.Range("D44").Value = "1.961.12"
$ws.Range("E44").Value = "  -3.45%  "
$ws.Range("E45").Value = "  -2.67%  "
$ws.Range("D46").Value = "'17.45"
$ws.Range("E46").Value = "  -7.01%  "
$ws.Range("D47").Value = "'9.60"
$ws.Range("E47").Value = "  -5.63%  "
$ws.Range("E48").Value = "  -4.93%  "
$ws.Range("D49").Value = "2.499.69"
$ws.Range("E49").Value = "  -2.23%  "
$ws.Range("D50").Value = "'52.13"
$ws.Range("E50").Value = "  -8.14%  "
$ws.Range("D51").Value = "'2.74"
$ws.Range("E51").Value = "  -5.16%  "